$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Headers
$ws.Range("W1").Value = "Utahgasprice_dollperMCF"
$ws.Range("X1").Value = "Utahcrudeprice_dollperBBL"

# Data values for Utah gas price ($/MCF) and Utah crude price ($/BBL)
$ws.Range("W2").Value = 5.214999999999999
$ws.Range("X2").Value = 84.702500000000001

$ws.Range("W3").Value = 5.8858333333333333
$ws.Range("X3").Value = 79.110000000000014

$ws.Range("W4").Value = 5.8908333333333331
$ws.Range("X4").Value = 40.118333333333332

$ws.Range("W5").Value = 5.4483333333333333
$ws.Range("X5").Value = 36.996666666666663

$ws.Range("W6").Value = 5.4858333333333347
$ws.Range("X6").Value = 44.310833333333328

$ws.Range("W7").Value = 5.21
$ws.Range("X7").Value = 56.95333333333334

$ws.Range("W8").Value = 4.9316666666666666
$ws.Range("X8").Value = 48.303333333333342

$ws.Range("W9").Value = 5.0183333333333318
$ws.Range("X9").Value = 34.637499999999996

$ws.Range("W10").Value = 5.4233333333333347
$ws.Range("X10").Value = 60.24916666666666

$ws.Range("W11").Value = 7.9741666666666662
$ws.Range("X11").Value = 81.140833333333333

$ws.Range("W12").Value = 9.9150000000000009
$ws.Range("X12").Value = 66.517499999999998

$ws.Range("W13").Value = 7.6291666666666664
$ws.Range("X13").Value = 62.831666666666671

# Comments with source links (Seth Lyman)
$excel.UserName = "Seth Lyman"

$c1 = $ws.Range("W1").AddComment("Seth Lyman:`nhttps://www.eia.gov/dnav/ng/hist/n3035ut3m.htm")
$c2 = $ws.Range("X1").AddComment("Seth Lyman:`nhttps://www.eia.gov/dnav/pet/hist/LeafHandler.ashx?n=PET&s=F004049__3&f=M")

# Freeze panes / view state to match new data extent
$ws.Range("G2").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true
$ws.Range("W2").Select() | Out-Null
